# Auto-generated script: apply scheduled market-data refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 416.66666
$ws.Range("I2").Value = 225
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 225
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -112
$ws.Range("N2").Value = -1026
$ws.Range("H26").Value = 49500
$ws.Range("J26").Value = 49000
$ws.Range("L26").Value = 49000
$ws.Range("N26").Value = -49688
$ws.Range("H98").Value = 560711.6
$ws.Range("I98").Value = 860433.1
$ws.Range("J98").Value = 4086
$ws.Range("K98").Value = 860433.1
$ws.Range("L98").Value = 4086
$ws.Range("M98").Value = -858935.1
$ws.Range("N98").Value = -7082
$ws.Range("H106").Value = 9343769
$ws.Range("I106").Value = 9343769
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 9343769
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -9343138
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 1000005
$ws.Range("I113").Value = 1000005
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000005
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -996751
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 560711.6
$ws.Range("I122").Value = 860433.1
$ws.Range("J122").Value = 4086
$ws.Range("K122").Value = 2581299.3
$ws.Range("L122").Value = 12258
$ws.Range("M122").Value = -2578849.3
$ws.Range("N122").Value = -17158
$ws.Range("H137").Value = 1898.8334
$ws.Range("I137").Value = 1800
$ws.Range("J137").Value = 1904.6471
$ws.Range("K137").Value = 5400
$ws.Range("L137").Value = 5713.9413
$ws.Range("M137").Value = -2850
$ws.Range("N137").Value = -10813.9413
$ws.Range("H138").Value = 24150904
$ws.Range("I138").Value = 1956554
$ws.Range("J138").Value = 83335830
$ws.Range("K138").Value = 5869662
$ws.Range("L138").Value = 250007490
$ws.Range("M138").Value = -5864522
$ws.Range("N138").Value = -250017770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1406999.8
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1406999.8
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1406999.8
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -1407299.8
$ws.Range("H32").Value = 3169.5186
$ws.Range("I32").Value = 2301.1702
$ws.Range("J32").Value = 8999.857
$ws.Range("K32").Value = 2301.1702
$ws.Range("L32").Value = 8999.857
$ws.Range("M32").Value = -2014.1702
$ws.Range("N32").Value = -9573.857
$ws.Range("H74").Value = 4424
$ws.Range("I74").Value = 1243.4584
$ws.Range("J74").Value = 10295.77
$ws.Range("K74").Value = 1243.4584
$ws.Range("L74").Value = 10295.77
$ws.Range("M74").Value = -369.4584
$ws.Range("N74").Value = -12043.77
$ws.Range("H77").Value = 4424
$ws.Range("I77").Value = 1243.4584
$ws.Range("J77").Value = 10295.77
$ws.Range("K77").Value = 6217.291999999999
$ws.Range("L77").Value = 51478.85000000001
$ws.Range("M77").Value = -1849.291999999999
$ws.Range("N77").Value = -60214.85000000001
$ws.Range("H97").Value = 11908738
$ws.Range("I97").Value = 15156480
$ws.Range("K97").Value = 15156480
$ws.Range("M97").Value = -15155984
$ws.Range("H132").Value = 1905.8474
$ws.Range("I132").Value = 1565.4166
$ws.Range("J132").Value = 3391.3635
$ws.Range("K132").Value = 4696.2498
$ws.Range("L132").Value = 10174.0905
$ws.Range("M132").Value = -2166.2498
$ws.Range("N132").Value = -15234.0905
$ws.Range("H139").Value = 68496.664
$ws.Range("J139").Value = 68496.664
$ws.Range("L139").Value = 68496.664
$ws.Range("N139").Value = -78776.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2181.4473
$ws.Range("I134").Value = 1298.8
$ws.Range("J134").Value = 5491.375
$ws.Range("K134").Value = 3896.4
$ws.Range("L134").Value = 16474.125
$ws.Range("M134").Value = -1361.4
$ws.Range("N134").Value = -21544.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2428.96
$ws.Range("I31").Value = 1498
$ws.Range("J31").Value = 3825.4
$ws.Range("K31").Value = 1498
$ws.Range("L31").Value = 3825.4
$ws.Range("M31").Value = -1203
$ws.Range("N31").Value = -4415.4
$ws.Range("H34").Value = 2428.96
$ws.Range("I34").Value = 1498
$ws.Range("J34").Value = 3825.4
$ws.Range("K34").Value = 1498
$ws.Range("L34").Value = 3825.4
$ws.Range("M34").Value = -1296
$ws.Range("N34").Value = -4229.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 135.12
$ws.Range("I38").Value = 179.81818
$ws.Range("K38").Value = 539.4545400000001
$ws.Range("M38").Value = -192.4545400000001
$ws.Range("H68").Value = 1092.5555
$ws.Range("I68").Value = 879.8226
$ws.Range("J68").Value = 1449.027
$ws.Range("K68").Value = 2639.4678
$ws.Range("L68").Value = 4347.081
$ws.Range("M68").Value = -1828.4678
$ws.Range("N68").Value = -5969.081
$ws.Range("H69").Value = 4127.5
$ws.Range("J69").Value = 4127.5
$ws.Range("L69").Value = 12382.5
$ws.Range("N69").Value = -14004.5
$ws.Range("H71").Value = 1092.5555
$ws.Range("I71").Value = 879.8226
$ws.Range("J71").Value = 1449.027
$ws.Range("K71").Value = 7918.403399999999
$ws.Range("L71").Value = 13041.243
$ws.Range("M71").Value = -3862.403399999999
$ws.Range("N71").Value = -21153.243
$ws.Range("H72").Value = 4127.5
$ws.Range("J72").Value = 4127.5
$ws.Range("L72").Value = 37147.5
$ws.Range("N72").Value = -45259.5
$ws.Range("H80").Value = 1129
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 1167.1666
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3501.4998
$ws.Range("M80").Value = -1764
$ws.Range("N80").Value = -5373.4998
$ws.Range("H83").Value = 1129
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 1167.1666
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 10504.4994
$ws.Range("M83").Value = -3420
$ws.Range("N83").Value = -19864.4994
$ws.Range("H131").Value = 3042.0657
$ws.Range("J131").Value = 3211.7017
$ws.Range("L131").Value = 9635.105100000001
$ws.Range("N131").Value = -19715.1051

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1098.2222
$ws.Range("I107").Value = 1098.2222
$ws.Range("K107").Value = 1098.2222
$ws.Range("M107").Value = 821.7778000000001
$ws.Range("H132").Value = 2608.244
$ws.Range("I132").Value = 2112.1724
$ws.Range("K132").Value = 6336.5172
$ws.Range("M132").Value = -3806.5172
$ws.Range("H138").Value = 65466.668
$ws.Range("J138").Value = 65466.668
$ws.Range("L138").Value = 65466.668
$ws.Range("N138").Value = -75746.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3000
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2774
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 868.6667
$ws.Range("I22").Value = 875.55554
$ws.Range("J22").Value = 848
$ws.Range("K22").Value = 875.55554
$ws.Range("L22").Value = 848
$ws.Range("M22").Value = -580.55554
$ws.Range("N22").Value = -1438
$ws.Range("H27").Value = 868.6667
$ws.Range("I27").Value = 875.55554
$ws.Range("J27").Value = 848
$ws.Range("K27").Value = 875.55554
$ws.Range("L27").Value = 848
$ws.Range("M27").Value = -768.55554
$ws.Range("N27").Value = -1062
$ws.Range("H104").Value = 27500
$ws.Range("J104").Value = 27500
$ws.Range("L104").Value = 27500
$ws.Range("N104").Value = -34488
$ws.Range("H135").Value = 41079.5
$ws.Range("J135").Value = 41079.5
$ws.Range("L135").Value = 41079.5
$ws.Range("N135").Value = -51219.5
$ws.Range("H136").Value = 2552.3865
$ws.Range("I136").Value = 1148.08
$ws.Range("J136").Value = 4400.1577
$ws.Range("K136").Value = 3444.24
$ws.Range("L136").Value = 13200.4731
$ws.Range("M136").Value = -894.2399999999998
$ws.Range("N136").Value = -18300.4731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 186372.67
$ws.Range("I75").Value = 500118
$ws.Range("K75").Value = 500118
$ws.Range("M75").Value = -499182
$ws.Range("H78").Value = 186372.67
$ws.Range("I78").Value = 500118
$ws.Range("K78").Value = 1500354
$ws.Range("M78").Value = -1495674
$ws.Range("H100").Value = 400.33334
$ws.Range("I100").Value = 200.66667
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 401.33334
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = 139.66666
$ws.Range("N100").Value = -2282
$ws.Range("H126").Value = 60227.707
$ws.Range("I126").Value = 78051.766
$ws.Range("J126").Value = 2299.5
$ws.Range("K126").Value = 234155.298
$ws.Range("L126").Value = 6898.5
$ws.Range("M126").Value = -231685.298
$ws.Range("N126").Value = -11838.5
$ws.Range("H132").Value = 11365590
$ws.Range("I132").Value = 17858832
$ws.Range("J132").Value = 2417.75
$ws.Range("K132").Value = 53576496
$ws.Range("L132").Value = 7253.25
$ws.Range("M132").Value = -53573966
$ws.Range("N132").Value = -12313.25
$ws.Range("H136").Value = 8800447
$ws.Range("I136").Value = 30394640
$ws.Range("J136").Value = 2812.5186
$ws.Range("K136").Value = 91183920
$ws.Range("L136").Value = 8437.5558
$ws.Range("M136").Value = -91181370
$ws.Range("N136").Value = -13537.5558
